# ECE251 assignment 3 v2 — applies the commit's edits to kemp_03.docx
# Strategy: each touched paragraph is replaced wholesale via
# Range.InsertXML() with the exact target run/proofErr/bookmark markup,
# which lets us reproduce multi-run splits (and proofErr spell/grammar
# tags) that a plain Find/Replace cannot produce.

$d = $word.ActiveDocument
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$enDash = [char]0x2013

function Set-ParagraphXml($paraIndex, $attrs, $inner) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $xml = "<w:p $w" + $attrs + ">" + $inner + "</w:p>"
    $rng.InsertXML($xml)
}

# Paragraph 2: "ECE251 Assignment 2 explanation" -> 3 runs (split the "2" -> "3")
Set-ParagraphXml 2 " w:rsidR='004A595D' w:rsidRDefault='004A595D'" (
    "<w:r><w:t xml:space='preserve'>ECE251 Assignment </w:t></w:r>" +
    "<w:r><w:t>3</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> explanation</w:t></w:r>"
)

# Paragraph 3: "-" -> long explanation about the eye-plot error
Set-ParagraphXml 3 " w:rsidR='00540803' w:rsidRDefault='00D36FAD' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>The error in the eye plot was because I was giving the filter a long-pulse train rather than an impulse train at the symbol frequency. This resulted in both the scaling issue and the eye-plot error. I was inadvertently causing ISI in my time domain </w:t></w:r>" +
    "<w:r><w:t>signal</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> because I was superimposing the impulse responses of the SRRC filter too close together by giving the filter a PAM </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>squarewave</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> instead of an impulse train.</w:t></w:r>"
)

# Paragraph 4: long "I noticed..." paragraph -> en dash placeholder "–"
Set-ParagraphXml 4 " w:rsidR='00E912F6' w:rsidRDefault='00D36FAD' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t>$enDash</w:t></w:r>"
)

# Paragraph 6: long "I'm noticing..." paragraph -> "Resolved error due to fix in part 1."
Set-ParagraphXml 6 " w:rsidR='00D36FAD' w:rsidRDefault='00D36FAD' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Resolved error </w:t></w:r>" +
    "<w:r><w:t>due to fix in part 1.</w:t></w:r>"
)

# Paragraph 7: en dash "–" placeholder -> plain hyphen "-"
Set-ParagraphXml 7 " w:rsidR='00D36FAD' w:rsidRDefault='00D36FAD' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t>-</w:t></w:r>"
)

# Paragraph 8: split the final run so the _GoBack bookmark moves in front of
# the trailing period ("...edge." -> "...edge" + bookmark + ".")
Set-ParagraphXml 8 " w:rsidR='00D36FAD' w:rsidRDefault='00D936CA' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t>I had to add an additional delay of .5*</w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>sps</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> to the start of the </w:t></w:r>" +
    "<w:r w:rsidR='00FE5CAD'><w:t xml:space='preserve'>sampling train to ensure that the center of each symbol </w:t></w:r>" +
    "<w:r w:rsidR='00EC386E'><w:t xml:space='preserve'>in y(t) </w:t></w:r>" +
    "<w:r w:rsidR='00FE5CAD'><w:t>is being sampled, rather than the edge</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "<w:r><w:t>.</w:t></w:r>"
)

# Paragraph 9: drop the now-relocated _GoBack bookmark from the document end
Set-ParagraphXml 9 " w:rsidR='00D936CA' w:rsidRDefault='00D936CA' w:rsidP='001644A3'" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
    "<w:r><w:t>I had to get rid of the last number of symbols equal to the span of the SRRC filter because they get pushed out due to the filter delay. Otherwise they will artificially drag down the bit-error-rate.</w:t></w:r>" +
    "<w:r w:rsidR='00FE5CAD'><w:t xml:space='preserve'> When I run the script with an alpha of 0.2, I usually get 2 to 4 errors. As far as I know, these errors are random and not corresponding to a systematic mistake.</w:t></w:r>" +
    "<w:r w:rsidR='00E467A4'><w:t xml:space='preserve'> When I run the script with an alpha of slightly greater than 0.2, (e.g. 0.25) these errors go away. 0.2 seems to be just on the threshold of where these errors arise.</w:t></w:r>"
)
